$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("run_1")
$ws.Range("F2").Value = 30.55284309387207
$ws.Range("F3").Value = 30.20093321800232
$ws.Range("F4").Value = 30.34845042228699
$ws.Range("F5").Value = 30.17138481140137
$ws.Range("F6").Value = 30.28632783889771
$ws.Range("F7").Value = 30.19257831573486
$ws.Range("F8").Value = 30.29971837997437
$ws.Range("F9").Value = 30.15856862068176
$ws.Range("F10").Value = 30.20747447013855
$ws.Range("F11").Value = 30.54825091362
$ws.Range("F12").Value = 30.23139953613281
$ws.Range("F13").Value = 29.75354814529419
$ws.Range("F14").Value = 29.7736713886261
$ws.Range("F15").Value = 29.64809727668763
$ws.Range("F16").Value = 31.31631135940552
$ws.Range("F17").Value = 31.40930461883545
$ws.Range("F18").Value = 30.63635301589965
$ws.Range("F19").Value = 30.10712838172913
$ws.Range("F20").Value = 30.26243305206299
$ws.Range("F21").Value = 30.42222666740417

$ws = $wb.Worksheets.Item("run_2")
$ws.Range("F2").Value = 30.50125789642334
$ws.Range("F3").Value = 30.1883978843689
$ws.Range("F4").Value = 30.17468476295471
$ws.Range("F5").Value = 30.22537612915039
$ws.Range("F6").Value = 30.42898106575012
$ws.Range("F7").Value = 30.36958050727844
$ws.Range("F8").Value = 30.34351801872253
$ws.Range("F9").Value = 30.50501132011414
$ws.Range("F10").Value = 30.383305311203
$ws.Range("F11").Value = 30.53819823265076
$ws.Range("F12").Value = 29.88113451004028
$ws.Range("F13").Value = 29.70088911056519
$ws.Range("F14").Value = 29.89678835868835
$ws.Range("F15").Value = 29.67514443397522
$ws.Range("F16").Value = 31.10500454902649
$ws.Range("F17").Value = 31.16188406944275
$ws.Range("F18").Value = 30.62839341163636
$ws.Range("F19").Value = 30.08226132392884
$ws.Range("F20").Value = 30.33508777618408
$ws.Range("F21").Value = 30.46755743026733

$ws = $wb.Worksheets.Item("run_3")
$ws.Range("F2").Value = 30.68663883209229
$ws.Range("F3").Value = 30.26375436782837
$ws.Range("F4").Value = 30.32671761512756
$ws.Range("F5").Value = 30.25673484802246
$ws.Range("F6").Value = 30.36309003829956
$ws.Range("F7").Value = 30.33768081665039
$ws.Range("F8").Value = 30.32880067825317
$ws.Range("F9").Value = 30.3119044303894
$ws.Range("F10").Value = 30.26532864570618
$ws.Range("F11").Value = 30.54590082168579
$ws.Range("F12").Value = 30.13988280296326
$ws.Range("F13").Value = 29.70302176475525
$ws.Range("F14").Value = 29.76651787757873
$ws.Range("F15").Value = 29.72830104827881
$ws.Range("F16").Value = 31.13151144981384
$ws.Range("F17").Value = 31.24899840354919
$ws.Range("F18").Value = 30.64975023269653
$ws.Range("F19").Value = 30.12477469444275
$ws.Range("F20").Value = 30.2613615989685
$ws.Range("F21").Value = 30.49258613586425

$ws = $wb.Worksheets.Item("run_4")
$ws.Range("F2").Value = 30.54592967033386
$ws.Range("F3").Value = 30.36660742759705
$ws.Range("F4").Value = 30.36813998222351
$ws.Range("F5").Value = 30.68802237510681
$ws.Range("F6").Value = 30.30759692192078
$ws.Range("F7").Value = 30.45486712455749
$ws.Range("F8").Value = 30.34749221801757
$ws.Range("F9").Value = 30.98499298095703
$ws.Range("F10").Value = 30.33475375175476
$ws.Range("F11").Value = 30.69737362861633
$ws.Range("F12").Value = 29.92249917984009
$ws.Range("F13").Value = 29.80490612983704
$ws.Range("F14").Value = 29.76012301445008
$ws.Range("F15").Value = 29.82798790931702
$ws.Range("F16").Value = 31.17224550247192
$ws.Range("F17").Value = 31.318852186203
$ws.Range("F18").Value = 30.60292482376098
$ws.Range("F19").Value = 30.46082544326782
$ws.Range("F20").Value = 30.48412680625916
$ws.Range("F21").Value = 30.6363615989685

$ws = $wb.Worksheets.Item("run_5")
$ws.Range("F2").Value = 30.51481199264526
$ws.Range("F3").Value = 30.36656355857849
$ws.Range("F4").Value = 30.27360200881958
$ws.Range("F5").Value = 30.39195013046265
$ws.Range("F6").Value = 30.42551565170288
$ws.Range("F7").Value = 30.48582410812378
$ws.Range("F8").Value = 30.34378623962402
$ws.Range("F9").Value = 30.42120933532715
$ws.Range("F10").Value = 30.2649827003479
$ws.Range("F11").Value = 30.75757908821106
$ws.Range("F12").Value = 29.87369394302368
$ws.Range("F13").Value = 29.83626770973206
$ws.Range("F14").Value = 29.87786245346069
$ws.Range("F15").Value = 29.8292543888092
$ws.Range("F16").Value = 31.04256510734558
$ws.Range("F17").Value = 31.25989723205566
$ws.Range("F18").Value = 30.58058476448059
$ws.Range("F19").Value = 30.41161751747132
$ws.Range("F20").Value = 30.38379859924316
$ws.Range("F21").Value = 30.67462277412415
